$d = $word.ActiveDocument

$pairs = @(
    @("98×34=", "34×78="),
    @("59×65=", "50×72="),
    @("28×55=", "33×41="),
    @("26×78=", "62×83="),
    @("90×48=", "84×89="),
    @("21×54=", "30×79="),
    @("18×43=", "85×89="),
    @("72×34=", "81×75="),
    @("31×11=", "21×30="),
    @("40×57=", "44×19="),
    @("48×98=", "87×93="),
    @("70×92=", "15×22="),
    @("64×63=", "57×14="),
    @("21×92=", "16×28="),
    @("45×80=", "50×25="),
    @("67×59=", "57×72="),
    @("35×79=", "86×22="),
    @("71×90=", "82×95="),
    @("71×60=", "25×84="),
    @("92×54=", "32×36="),
    @("97×23=", "33×99="),
    @("79×98=", "98×65="),
    @("23×85=", "22×19="),
    @("36×53=", "88×56="),
    @("79×13=", "91×97=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
